# "Customer misbehavior" meta-analysis workbook update:
#  - study_data: delete the unused, empty column G (shifts H:X left to G:W),
#    which also moves the review comment from S3 to R3.
#  - study_data becomes the active/selected sheet (was es_data), with a new
#    selection at F7 and the frozen/top-left view reset.
#  - es_data: gets an AutoFilter over its header row, loses its "selected
#    tab" flag, and the workbook gains the matching hidden _FilterDatabase
#    defined name that Excel creates when AutoFilter is turned on.

$wb = $excel.ActiveWorkbook

$wsStudy = $wb.Worksheets.Item("study_data")
$wsEs = $wb.Worksheets.Item("es_data")

# --- study_data: remove the empty column G -------------------------------
$old = $wsStudy.Range("S3").Comment
$commentText = $old.Text()
$old.Delete()

$hyperlinkUrl = "https://sci-hub.wf/10.1177/0093854810374282"
$wsStudy.Range("X2").Hyperlinks.Delete()

$wsStudy.Columns.Item(7).Delete()

$newComment = $wsStudy.Range("R3").AddComment($commentText)
$wsStudy.Hyperlinks.Add($wsStudy.Range("W2"), $hyperlinkUrl) | Out-Null

# --- es_data: add AutoFilter + the matching defined name -----------------
$wsEs.Range("A1:H1").AutoFilter()
$filterName = $wsEs.Names.Add("_xlnm._FilterDatabase", "=es_data!`$A`$1:`$H`$1")
$filterName.Visible = $false

# --- view/selection updates ------------------------------------------------
$wsEs.Range("C8").Select()

$wsStudy.Activate()
$wsStudy.Application.ActiveWindow.ScrollColumn = 1
$wsStudy.Application.ActiveWindow.ScrollRow = 1
$wsStudy.Range("F7").Select()
